# AFDP-157 - Add Access Control List to SOLR documents - apply assignment
# and data access control rules to case files.
#
# This extends the "Access Control Rules" rule table on Sheet1 with five
# new rows (23-27) describing Complaint lockout and Case File rules, shifts
# the previous filler rows down, and appends one more blank templated row
# (row 32) at the bottom of the table to keep the same amount of spare rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 23 - Complaint - Lockout No Access Users
# ---------------------------------------------------------------------
$ws.Rows.Item(23).RowHeight = 23.5
$ws.Range("B23").Value = "Complaint – Lockout No Access Users"
$ws.Range("C23").Value = "COMPLAINT"
$ws.Range("G23").Value = "mandatory deny read to No Access"

# ---------------------------------------------------------------------
# Row 24 - Case File - Assignee Read Access
# ---------------------------------------------------------------------
$ws.Rows.Item(24).RowHeight = 23.65
$ws.Range("B24").Value = "Case File – Assignee Read Access"
$ws.Range("C24").Value = "CASE_FILE"
$ws.Range("C24").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("G24").Value = "grant read to assignee"

# ---------------------------------------------------------------------
# Row 25 - Case File - Restrict Access to Drafts
# ---------------------------------------------------------------------
$ws.Rows.Item(25).RowHeight = 23.65
$ws.Range("B25").Value = "Case File – Restrict Access to Drafts"
$ws.Range("C25").Value = "CASE_FILE"
$ws.Range("C25").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("D25").Value = "status == 'DRAFT'"
$ws.Range("D25").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("G25").Value = "deny read to *"

# ---------------------------------------------------------------------
# Row 26 - Case File - Grant Access to non-Drafts
# ---------------------------------------------------------------------
$ws.Rows.Item(26).RowHeight = 23.65
$ws.Range("B26").Value = "Case File – Grant Access to non-Drafts"
$ws.Range("C26").Value = "CASE_FILE"
$ws.Range("C26").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("D26").Value = "status != 'DRAFT'"
$ws.Range("D26").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("G26").Value = "grant read to *"

# ---------------------------------------------------------------------
# Row 27 - Case File - Lockout No Access Users
# ---------------------------------------------------------------------
$ws.Rows.Item(27).RowHeight = 23.5
$ws.Range("B27").Value = "Case File – Lockout No Access Users"
$ws.Range("C27").Value = "CASE_FILE"
$ws.Range("G27").Value = "mandatory deny read to No Access"

# ---------------------------------------------------------------------
# Row 32 - extend the table with one more blank templated row, matching
# the look (borders/shading) of the other blank rows above it (20-31).
# ---------------------------------------------------------------------
$ws.Rows.Item(32).Insert()
$ws.Range("B31").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("C31:G31").Copy()
$ws.Range("C32:G32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Restore the active selection to where the editor ended up after the
# change (cell D29).
# ---------------------------------------------------------------------
$ws.Range("D29").Select() | Out-Null
